$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A20").Value = "telecomm/gsm/runme_large.sh"
$ws.Range("B20").Value = 0.15
$ws.Range("C20").Value = 0.14
$ws.Range("D20").Value = 0
$ws.Range("B21").Select()
